# Price update for 2026-02-07 - append a new tracking row at the bottom
# of the sheet (row 38): Date, Price, Discount, Incredible.
#
# The new values are written as formulas first (="text") so Excel's
# automatic "looks like a date/number" input-parsing doesn't kick in and
# reformat the cell (which would otherwise allocate a new cell style).
# They are then copied and pasted back as values in place, which converts
# the formula cells into plain shared-string literals matching the rest
# of the sheet, while leaving styles.xml untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Range("A$row").Formula = '="2026-02-07"'
$ws.Range("B$row").Formula = '="3252100"'
$ws.Range("C$row").Formula = '="6"'
$ws.Range("D$row").Formula = '="0"'

$target = $ws.Range("A$row" + ":D$row")
$target.Copy()
$target.PasteSpecial(-4163)
